# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型"
# sheets to reflect the latest data scrape, per commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F ("想去人数")
$updates = @{
    2  = 3428
    6  = 1671
    7  = 483
    8  = 388
    13 = 239
    15 = 58
    19 = 31
    22 = 132
    25 = 299
    27 = 51
    30 = 596
    31 = 2399
    35 = 692
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
